$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

# New row 2 values (newest data point inserted at the top of the rolling window)
$ws.Range("B2").Value = 0.1660323612164193
$ws.Range("C2").Value = 0.468872571632237
$ws.Range("D2").Value = 0.3770675980855352
$ws.Range("E2").Value = 0.6140583018619121
$ws.Range("F2").Value = 0.6119356173368786
$ws.Range("G2").Value = 15

# Row 3 (previously row 2's values)
$ws.Range("B3").Value = 0.2261183143199825
$ws.Range("C3").Value = 0.4923620343352079
$ws.Range("D3").Value = 0.4821095245989872
$ws.Range("E3").Value = 0.6943410722397079
$ws.Range("F3").Value = 0.6812725906144649
$ws.Range("G3").Value = 14

# Row 4 (previously row 3's values)
$ws.Range("B4").Value = 0.2360243239862733
$ws.Range("C4").Value = 0.5248832413428339
$ws.Range("D4").Value = 0.37481815022521
$ws.Range("E4").Value = 0.6122239379714011
$ws.Range("F4").Value = 0.5879653258807318
$ws.Range("G4").Value = 13

# Row 5 (previously row 4's values)
$ws.Range("B5").Value = 0.3051580102940423
$ws.Range("C5").Value = 0.3983834828835194
$ws.Range("D5").Value = 0.2461501466550942
$ws.Range("E5").Value = 0.4961352100537657
$ws.Range("F5").Value = 0.4085834536877719
$ws.Range("G5").Value = 12

# Row 6 (previously row 5's values)
$ws.Range("B6").Value = 0.3404744896530312
$ws.Range("C6").Value = 0.5484269867923621
$ws.Range("D6").Value = 0.5429546666577961
$ws.Range("E6").Value = 0.736854576329547
$ws.Range("F6").Value = 0.6853721378992834
$ws.Range("G6").Value = 11

# Row 7 (previously row 6's values)
$ws.Range("B7").Value = 0.3329151699208477
$ws.Range("C7").Value = 0.5578246582442683
$ws.Range("D7").Value = 0.4732217363555192
$ws.Range("E7").Value = 0.6879111398687473
$ws.Range("F7").Value = 0.6345507824806217
$ws.Range("G7").Value = 10

# Row 8 (previously row 7's values)
$ws.Range("B8").Value = 0.2814885994813455
$ws.Range("C8").Value = 0.5772087615654313
$ws.Range("D8").Value = 0.5238078509972534
$ws.Range("E8").Value = 0.7237457087936711
$ws.Range("F8").Value = 0.7072082591282391
$ws.Range("G8").Value = 9

# Row 9 (previously row 8's values)
$ws.Range("B9").Value = 0.3213393071819032
$ws.Range("C9").Value = 0.6123030108606557
$ws.Range("D9").Value = 0.5290292607895575
$ws.Range("E9").Value = 0.7273439769390804
$ws.Range("F9").Value = 0.6975633594976256
$ws.Range("G9").Value = 8

# Row 10 (previously row 9's values)
$ws.Range("B10").Value = 0.4424481932715923
$ws.Range("C10").Value = 0.4424481932715923
$ws.Range("D10").Value = 0.2455568012204569
$ws.Range("E10").Value = 0.4955368817963572
$ws.Range("F10").Value = 0.2410306973665816
$ws.Range("G10").Value = 7

# Row 11 (previously row 10's values)
$ws.Range("B11").Value = 0.3495035403546135
$ws.Range("C11").Value = 0.4324022405235847
$ws.Range("D11").Value = 0.3358454641318753
$ws.Range("E11").Value = 0.5795217546666176
$ws.Range("F11").Value = 0.5063904494495919
$ws.Range("G11").Value = 6
